$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "libraryPreparer" column (B) was filled with "BROWN" for every sample
# row; normalize it to "H.BROWN" across the whole data range (rows 2-53).
$ws.Range("B2:B53").Value = "H.BROWN"

# Un-hide the previously hidden metadata columns D:J (keep their width).
$ws.Range("D1:J1").EntireColumn.Hidden = $false

# Restore the view to the top of the sheet with F14 selected.
$ws.Range("F14").Select()
